$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")
$wsTot = $wb.Worksheets.Item("Totaux")

# --- Update existing hour values (Journal sheet) ---
$ws.Range("C7").Value = 1/24
$ws.Range("C16").Value = 2.25/24
$ws.Range("C18").Value = 1.25/24

# --- Re-purpose the description of row 18 ---
$ws.Range("E18").Value = "Rédaction du rapport de projet"

# --- Add a new row to the Journal table ---
$tbl = $ws.ListObjects.Item("Tableau1")
$newRow = $tbl.ListRows.Add()

$ws.Range("A19").Value = 44977
$ws.Range("A19").NumberFormat = "d-mmm"
$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("A19").VerticalAlignment = -4108

$ws.Range("B19").Value = 3
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("B19").VerticalAlignment = -4108

$ws.Range("C19").Value = 0.5/24
$ws.Range("C19").NumberFormat = "hh/mm"" h"";@"
$ws.Range("C19").HorizontalAlignment = -4108
$ws.Range("C19").VerticalAlignment = -4108

$ws.Range("D19").Value = "Analyse"
$ws.Range("D19").HorizontalAlignment = -4108
$ws.Range("D19").VerticalAlignment = -4108

$ws.Range("E19").Value = 'Analyse des templates dans la partie "gestion du parc"'

# --- Totaux sheet: add the missing weekly sum for the new week-3 rows ---
$wsTot.Range("B7").Formula = "=SUM(Journal!C16:C19)"

# --- Restore the selections recorded in the workbook ---
$ws.Range("D30").Select()
$wsTot.Range("A22").Select()
$ws.Activate()
